# Hortaliza, Vega Modelo de Temuco - Coliflor: add a new weekly price report row.
# A new record (week of 2021-10-22) is inserted at row 218, pushing the
# existing rows 218-252 down to 219-253; row 1 (header) and rows 2-217 are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 218, shifting rows 218:252 down to 219:253.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(218, 1).Value  = 10
$ws.Cells.Item(218, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(218, 3).Value  = "La Araucanía"
$ws.Cells.Item(218, 4).Value  = 44491
$ws.Cells.Item(218, 5).Value  = 9
$ws.Cells.Item(218, 6).Value  = 100112008
$ws.Cells.Item(218, 7).Value  = "Coliflor"
$ws.Cells.Item(218, 8).Value  = "Sin especificar"
$ws.Cells.Item(218, 9).Value  = "Primera"
$ws.Cells.Item(218, 10).Value = 3200
$ws.Cells.Item(218, 11).Value = 800
$ws.Cells.Item(218, 12).Value = 900
$ws.Cells.Item(218, 13).Value = 839
$ws.Cells.Item(218, 14).Value = "$/unidad"
$ws.Cells.Item(218, 15).Value = "Región Metropolitana"
$ws.Cells.Item(218, 16).Value = 839
$ws.Cells.Item(218, 17).Value = 1
$ws.Cells.Item(218, 18).Value = "Hortaliza"
